# smartart-picture-strip.pptx: fix lack of spacing around the picture list.
#
# The SmartArt graphic frame is moved/resized (the snake/spacing algorithm
# now reserves more room around the picture list) and a new title textbox
# ("Foo Bar Baz Blah") is added above it so the two no longer overlap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Reposition + resize the SmartArt graphic frame -----------------
# Target EMU: off (2419200, 835200), ext (6096000, 5418000).
# Shape.Left/Top/Width/Height are in points (1 pt = 12700 EMU); the literals
# below are float32-exact so they round-trip to precisely those EMU values.
$gf = $s.Shapes.Item(1)
$gf.Left = 190.48818969726562
$gf.Top = 65.76378631591797
$gf.Width = 480.0
$gf.Height = 426.61419677734375

# --- 2. Add the new title textbox above the SmartArt --------------------
# Target EMU: off (457200, 273600), ext (8229600, 1144800).
$tb = $s.Shapes.AddTextbox(1, 36.0, 21.54330825805664, 648.0, 90.14173889160156)

$tr = $tb.TextFrame.TextRange
$tr.Text = "Foo Bar "
$tr.Font.Size = 44
$tr.Font.Name = "+mj-lt"
[void]$tr.InsertAfter("Baz Blah")

$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

# Re-assert the exact geometry (autofit may have touched the height while
# the text/font were being applied above).
$tb.Left = 36.0
$tb.Top = 21.54330825805664
$tb.Width = 648.0
$tb.Height = 90.14173889160156
